# Apply the "updated data and reset archive" edit to the selfemployed_tosearch sheet.
# A new record (Brian McKeon / Department of State) was inserted at row 25, which shifts
# every subsequent Department of State / Veterans Affairs / EPA / EOP / FCC / Intelligence
# Community / NASA person record down by one row (through row 43); the trailing
# National Security Council rows (43-44 originally) collapse back to the original two rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 25 -> "Brian McKeon"
$ws.Range("B25").Value = 'Brian McKeon'
$ws.Range("H25").Value = 'Transition — PT Fund, Inc.'
$ws.Range("I25").Value = 'Brian McKeonSelf-employedDepartment of State'
$ws.Range("J25").Value = 'Brian McKeonSelf-employed'

# Row 26 -> "Fouad Saad"
$ws.Range("B26").Value = 'Fouad Saad'
$ws.Range("D26").Value = 'Y'
$ws.Range("I26").Value = 'Fouad SaadSelf-employedDepartment of State'
$ws.Range("J26").Value = 'Fouad SaadSelf-employed'

# Row 27 -> "Gentry Smith"
$ws.Range("B27").Value = 'Gentry Smith'
$ws.Range("I27").Value = 'Gentry SmithSelf-employedDepartment of State'
$ws.Range("J27").Value = 'Gentry SmithSelf-employed'

# Row 28 -> "Puneet Talwar"
$ws.Range("A28").Value = 'Department of State'
$ws.Range("B28").Value = 'Puneet Talwar'
$ws.Range("I28").Value = 'Puneet TalwarSelf-employedDepartment of State'
$ws.Range("J28").Value = 'Puneet TalwarSelf-employed'

# Row 29 -> "Chris Diaz"
$ws.Range("B29").Value = 'Chris Diaz'
$ws.Range("I29").Value = 'Chris DiazSelf-employedDepartment of Veterans Affairs'
$ws.Range("J29").Value = 'Chris DiazSelf-employed'

# Row 30 -> "Maryanne Donaghy"
$ws.Range("A30").Value = 'Department of Veterans Affairs'
$ws.Range("B30").Value = 'Maryanne Donaghy'
$ws.Range("I30").Value = 'Maryanne DonaghySelf-employedDepartment of Veterans Affairs'
$ws.Range("J30").Value = 'Maryanne DonaghySelf-employed'

# Row 31 -> "Michael McCabe"
$ws.Range("B31").Value = 'Michael McCabe'
$ws.Range("I31").Value = 'Michael McCabeSelf-employedEnvironmental Protection Agency'
$ws.Range("J31").Value = 'Michael McCabeSelf-employed'

# Row 32 -> "Luseni Pieh"
$ws.Range("A32").Value = 'Environmental Protection Agency'
$ws.Range("B32").Value = 'Luseni Pieh'
$ws.Range("H32").Value = 'Volunteer'
$ws.Range("I32").Value = 'Luseni PiehSelf-employedEnvironmental Protection Agency'
$ws.Range("J32").Value = 'Luseni PiehSelf-employed'

# Row 33 -> "Christian Peele"
$ws.Range("A33").Value = 'Executive Office of the President, Management and Administration'
$ws.Range("B33").Value = 'Christian Peele'
$ws.Range("H33").Value = 'Transition — PT Fund, Inc.'
$ws.Range("I33").Value = 'Christian PeeleSelf-employedExecutive Office of the President, Management and Administration'
$ws.Range("J33").Value = 'Christian PeeleSelf-employed'

# Row 34 -> "Mignon Clyburn"
$ws.Range("B34").Value = 'Mignon Clyburn'
$ws.Range("I34").Value = 'Mignon ClyburnSelf-employedFederal Communications Commission'
$ws.Range("J34").Value = 'Mignon ClyburnSelf-employed'

# Row 35 -> "Paul de Sa"
$ws.Range("A35").Value = 'Federal Communications Commission'
$ws.Range("B35").Value = 'Paul de Sa'
$ws.Range("C35").ClearContents()
$ws.Range("I35").Value = 'Paul de SaSelf-employedFederal Communications Commission'
$ws.Range("J35").Value = 'Paul de SaSelf-employed'

# Row 36 -> "Stephanie O’Sullivan"
$ws.Range("B36").Value = 'Stephanie O’Sullivan'
$ws.Range("C36").Value = 'Y'
$ws.Range("I36").Value = 'Stephanie O’Sullivan, Team LeadSelf-employedIntelligence Community'
$ws.Range("J36").Value = 'Stephanie O’Sullivan, Team LeadSelf-employed'

# Row 37 -> "Harry Coker"
$ws.Range("B37").Value = 'Harry Coker'
$ws.Range("I37").Value = 'Harry CokerSelf-employedIntelligence Community'
$ws.Range("J37").Value = 'Harry CokerSelf-employed'

# Row 38 -> "Dawn Eilenberger"
$ws.Range("B38").Value = 'Dawn Eilenberger'
$ws.Range("I38").Value = 'Dawn EilenbergerSelf-employedIntelligence Community'
$ws.Range("J38").Value = 'Dawn EilenbergerSelf-employed'

# Row 39 -> "Justin Jackson"
$ws.Range("B39").Value = 'Justin Jackson'
$ws.Range("I39").Value = 'Justin JacksonSelf-employedIntelligence Community'
$ws.Range("J39").Value = 'Justin JacksonSelf-employed'

# Row 40 -> "Ron Moultrie"
$ws.Range("B40").Value = 'Ron Moultrie'
$ws.Range("I40").Value = 'Ron MoultrieSelf-employedIntelligence Community'
$ws.Range("J40").Value = 'Ron MoultrieSelf-employed'

# Row 41 -> "Shaun Murphy"
$ws.Range("B41").Value = 'Shaun Murphy'
$ws.Range("I41").Value = 'Shaun MurphySelf-employedIntelligence Community'
$ws.Range("J41").Value = 'Shaun MurphySelf-employed'

# Row 42 -> "Bruce Pease"
$ws.Range("A42").Value = 'Intelligence Community'
$ws.Range("B42").Value = 'Bruce Pease'
$ws.Range("I42").Value = 'Bruce PeaseSelf-employedIntelligence Community'
$ws.Range("J42").Value = 'Bruce PeaseSelf-employed'

# Row 43 -> "Pam Melroy"
$ws.Range("A43").Value = 'National Aeronautics and Space Administration'
$ws.Range("B43").Value = 'Pam Melroy'
$ws.Range("H43").Value = 'Volunteer'
$ws.Range("I43").Value = 'Pam MelroySelf-employedNational Aeronautics and Space Administration'
$ws.Range("J43").Value = 'Pam MelroySelf-employed'
